$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.006795699771644
$ws.Range("D2").Value = 1.026513451475777
$ws.Range("E2").Value = 1.009532559647175
$ws.Range("F2").Value = 1.00497242631925
$ws.Range("I2").Value = 1.028140619206577
$ws.Range("J2").Value = 1.012071080893318
$ws.Range("K2").Value = 1.029335878222202
$ws.Range("L2").Value = 1.012405490914704
$ws.Range("M2").Value = 1.007859237321635
$ws.Range("N2").Value = 1.013508337914362

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.008197251359901
$ws.Range("D3").Value = 1.026832767717716
$ws.Range("E3").Value = 1.010735825058346
$ws.Range("F3").Value = 1.007032313929777
$ws.Range("I3").Value = 1.028119731977348
$ws.Range("J3").Value = 1.013101395776095
$ws.Range("K3").Value = 1.029464047931005
$ws.Range("L3").Value = 1.01341151132103
$ws.Range("M3").Value = 1.009718434850895
$ws.Range("N3").Value = 1.014540115962452

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.009102639152691
$ws.Range("D4").Value = 1.027038990663314
$ws.Range("E4").Value = 1.011513358573906
$ws.Range("F4").Value = 1.008363069200797
$ws.Range("I4").Value = 1.02810454699596
$ws.Range("J4").Value = 1.013766239260105
$ws.Range("K4").Value = 1.029545830079758
$ws.Range("L4").Value = 1.014060866618908
$ws.Range("M4").Value = 1.010918994503388
$ws.Range("N4").Value = 1.015205903600468

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.009482911626778
$ws.Range("D5").Value = 1.027125589626558
$ws.Range("E5").Value = 1.011839986557137
$ws.Range("F5").Value = 1.008922028095641
$ws.Range("I5").Value = 1.028097762610936
$ws.Range("J5").Value = 1.014045305955706
$ws.Range("K5").Value = 1.029579933844576
$ws.Range("L5").Value = 1.014333476974628
$ws.Range("M5").Value = 1.011423136912126
$ws.Range("N5").Value = 1.015485366602786

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.009546740631775
$ws.Range("D6").Value = 1.027140124207022
$ws.Range("E6").Value = 1.011894814532448
$ws.Range("F6").Value = 1.009015851539715
$ws.Range("I6").Value = 1.028096599963333
$ws.Range("J6").Value = 1.014092137257115
$ws.Range("K6").Value = 1.029585643700015
$ws.Range("L6").Value = 1.014379227428651
$ws.Range("M6").Value = 1.011507751549646
$ws.Range("N6").Value = 1.015532264410014

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.009107721746463
$ws.Range("D7").Value = 1.027040148187686
$ws.Range("E7").Value = 1.011517723955896
$ws.Range("F7").Value = 1.008370539941099
$ws.Range("I7").Value = 1.028104457918005
$ws.Range("J7").Value = 1.013769969858014
$ws.Range("K7").Value = 1.029546286868166
$ws.Range("L7").Value = 1.014064510731977
$ws.Range("M7").Value = 1.010925733108616
$ws.Range("N7").Value = 1.015209639496253

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.007269678681445
$ws.Range("D8").Value = 1.026621446480098
$ws.Range("E8").Value = 1.009939432213309
$ws.Range("F8").Value = 1.005669027589001
$ws.Range("I8").Value = 1.02813390540539
$ws.Range("J8").Value = 1.012419665540274
$ws.Range("K8").Value = 1.029379431388346
$ws.Range("L8").Value = 1.012745816713894
$ws.Range("M8").Value = 1.008488082817075
$ws.Range("N8").Value = 1.013857417591499

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.004018848159921
$ws.Range("D9").Value = 1.025880711243438
$ws.Range("E9").Value = 1.007149868128996
$ws.Range("F9").Value = 1.000891439981142
$ws.Range("I9").Value = 1.028173049635261
$ws.Range("J9").Value = 1.010025865257892
$ws.Range("K9").Value = 1.029076656696341
$ws.Range("L9").Value = 1.010409523095615
$ws.Range("M9").Value = 1.00417296713567
$ws.Range("N9").Value = 1.011460217838144

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.001843021329063
$ws.Range("D10").Value = 1.025385058414673
$ws.Range("E10").Value = 1.005284085024044
$ws.Range("F10").Value = 0.9976936178203025
$ws.Range("I10").Value = 1.028190631074268
$ws.Range("J10").Value = 1.008419891959954
$ws.Range("K10").Value = 1.02886901662911
$ws.Range("L10").Value = 1.008843141418668
$ws.Range("M10").Value = 1.001281905842057
$ws.Range("N10").Value = 1.009851963873915

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.000898696951886
$ws.Range("D11").Value = 1.025170031100778
$ws.Range("E11").Value = 1.004474645911796
$ws.Range("F11").Value = 0.9963056270408939
$ws.Range("I11").Value = 1.028196236062541
$ws.Range("J11").Value = 1.00772199414136
$ws.Range("K11").Value = 1.028777753015332
$ws.Range("L11").Value = 1.008162693502325
$ws.Range("M11").Value = 1.000026403695264
$ws.Range("N11").Value = 1.009153074960369

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.000547594827383
$ws.Range("D12").Value = 1.02509010175386
$ws.Range("E12").Value = 1.004173744240184
$ws.Range("F12").Value = 0.9957895443653946
$ws.Range("I12").Value = 1.028198017265132
$ws.Range("J12").Value = 1.007462379703838
$ws.Range("K12").Value = 1.028743651841878
$ws.Range("L12").Value = 1.007909607544663
$ws.Range("M12").Value = 0.9995594842477433
$ws.Range("N12").Value = 1.008893091840568

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.000622922897103
$ws.Range("D13").Value = 1.025107249484421
$ws.Range("E13").Value = 1.004238299660311
$ws.Range("F13").Value = 0.9959002698105656
$ws.Range("I13").Value = 1.028197648785465
$ws.Range("J13").Value = 1.007518085399739
$ws.Range("K13").Value = 1.028750975765278
$ws.Range("L13").Value = 1.00796391072889
$ws.Range("M13").Value = 0.9996596661921244
$ws.Range("N13").Value = 1.008948876644946

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.000869681699599
$ws.Range("D14").Value = 1.025163425298077
$ws.Range("E14").Value = 1.004449778228142
$ws.Range("F14").Value = 0.9962629782525639
$ws.Range("I14").Value = 1.028196389426981
$ws.Range("J14").Value = 1.007700542215689
$ws.Range("K14").Value = 1.028774938309027
$ws.Range("L14").Value = 1.008141780263487
$ws.Range("M14").Value = 0.9999878197400796
$ws.Range("N14").Value = 1.009131592570504

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.001021672823927
$ws.Range("D15").Value = 1.025198029374264
$ws.Range("E15").Value = 1.004580045112463
$ws.Range("F15").Value = 0.9964863851390772
$ws.Range("I15").Value = 1.028195573671256
$ws.Range("J15").Value = 1.007812908730181
$ws.Range("K15").Value = 1.028789675730796
$ws.Range("L15").Value = 1.00825132662824
$ws.Range("M15").Value = 1.000189929797714
$ws.Range("N15").Value = 1.009244118658335

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.001905646230184
$ws.Range("D16").Value = 1.025399320714895
$ws.Range("E16").Value = 1.005337771633596
$ws.Range("F16").Value = 0.9977856624107103
$ws.Range("I16").Value = 1.028190216864643
$ws.Range("J16").Value = 1.008466155736214
$ws.Range("K16").Value = 1.028875045105079
$ws.Range("L16").Value = 1.008888253624534
$ws.Range("M16").Value = 1.00136515060923
$ws.Range("N16").Value = 1.009898293350044

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.002459548863428
$ws.Range("D17").Value = 1.025525478374841
$ws.Range("E17").Value = 1.005812654693122
$ws.Range("F17").Value = 0.99859976137455
$ws.Range("I17").Value = 1.028186319765488
$ws.Range("J17").Value = 1.00887524489243
$ws.Range("K17").Value = 1.02892823368536
$ws.Range("L17").Value = 1.00928718825863
$ws.Range("M17").Value = 1.002101342565271
$ws.Range("N17").Value = 1.010307963459785

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.002782421593638
$ws.Range("D18").Value = 1.025599024633332
$ws.Range("E18").Value = 1.006089497654567
$ws.Range("F18").Value = 0.9990742934944372
$ws.Range("I18").Value = 1.02818385285959
$ws.Range("J18").Value = 1.009113618959025
$ws.Range("K18").Value = 1.02895912692504
$ws.Range("L18").Value = 1.009519668882094
$ws.Range("M18").Value = 1.002530399804979
$ws.Range("N18").Value = 1.010546676044896

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.002892477678187
$ws.Range("D19").Value = 1.025624095254975
$ws.Range("E19").Value = 1.00618386904272
$ws.Range("F19").Value = 0.9992360434449757
$ws.Range("I19").Value = 1.028182978809425
$ws.Range("J19").Value = 1.009194857778418
$ws.Range("K19").Value = 1.028969638492465
$ws.Range("L19").Value = 1.009598903173331
$ws.Range("M19").Value = 1.002676638324288
$ws.Range("N19").Value = 1.010628030232731

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.002400142081356
$ws.Range("D20").Value = 1.025511946906639
$ws.Range("E20").Value = 1.005761719632777
$ws.Range("F20").Value = 0.9985124493247259
$ws.Range("I20").Value = 1.028186757925705
$ws.Range("J20").Value = 1.008831378456744
$ws.Range("K20").Value = 1.028922540567341
$ws.Range("L20").Value = 1.009244408248378
$ws.Range("M20").Value = 1.002022392601343
$ws.Range("N20").Value = 1.010264034728729

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.000797026801523
$ws.Range("D21").Value = 1.025146884520153
$ws.Range("E21").Value = 1.004387509750594
$ws.Range("F21").Value = 0.9961561842308329
$ws.Range("I21").Value = 1.028196768571165
$ws.Range("J21").Value = 1.007646823902336
$ws.Range("K21").Value = 1.028767887494854
$ws.Range("L21").Value = 1.00808941149181
$ws.Range("M21").Value = 0.9998912026063483
$ws.Range("N21").Value = 1.009077797970983

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 0.9997871235086685
$ws.Range("D22").Value = 1.024917017289605
$ws.Range("E22").Value = 1.003522097095154
$ws.Range("F22").Value = 0.9946716801406038
$ws.Range("I22").Value = 1.028201323127522
$ws.Range("J22").Value = 1.006899820848372
$ws.Range("K22").Value = 1.028669484195791
$ws.Range("L22").Value = 1.007361263806506
$ws.Range("M22").Value = 0.9985479300973373
$ws.Range("N22").Value = 1.008329734087001

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.000322681920053
$ws.Range("D23").Value = 1.025038905493158
$ws.Range("E23").Value = 1.003981003448065
$ws.Range("F23").Value = 0.9954589386253812
$ws.Range("I23").Value = 1.028199073251875
$ws.Range("J23").Value = 1.007296035183665
$ws.Range("K23").Value = 1.028721759700692
$ws.Range("L23").Value = 1.007747456589261
$ws.Range("M23").Value = 0.9992603450241445
$ws.Range("N23").Value = 1.008726511092096

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.002426986095816
$ws.Range("D24").Value = 1.025518061316701
$ws.Range("E24").Value = 1.005784735452404
$ws.Range("F24").Value = 0.9985519028639323
$ws.Range("I24").Value = 1.028186560538906
$ws.Range("J24").Value = 1.008851200554655
$ws.Range("K24").Value = 1.028925113446604
$ws.Range("L24").Value = 1.009263739346227
$ws.Range("M24").Value = 1.00205806777604
$ws.Range("N24").Value = 1.010283884976293

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.004860741150578
$ws.Range("D25").Value = 1.026072541693807
$ws.Range("E25").Value = 1.007872076959016
$ws.Range("F25").Value = 1.002128720824313
$ws.Range("I25").Value = 1.028164434695646
$ws.Range("J25").Value = 1.010646469523226
$ws.Range("K25").Value = 1.029155958733625
$ws.Range("L25").Value = 1.011015042468879
$ws.Range("M25").Value = 1.005290970498297
$ws.Range("N25").Value = 1.012081703432719
